# "changed counts to 0"
#
# Table row "Total_Citations_(2023" previously summarised as
#   "8,828 (range: 11 to 169,257)"
# is updated to reflect a minimum of 0 (and a new mean of 7,170):
#   "7,170 (range: 0 to 169,257)"
#
# Because the minimum became 0, the separate "Unknown" row directly below
# it (which held the count of missing/unknown values, previously 86) is
# no longer applicable and is removed from the table - every row beneath
# it shifts up by one, and the trailing merged footer cell moves from
# A84:B84 to A83:B83 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the summary text for "Total_Citations_(2023" (row 79, column B).
$ws.Range("B79").Value2 = "7,170 (range: 0 to 169,257)"

# Remove the now-obsolete "    Unknown" / 86 row entirely; Excel shifts
# every subsequent row up by one and keeps the merged footer cell in sync.
$ws.Rows("80:80").Delete()
